$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.469.50"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.677.49"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'217.36"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").Value = "'0.5309"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.2695"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("E10").Value = "  +5.46%  "
$ws.Range("D11").Value = "'0.07820"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "1.682.66"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "'0.5583"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "0.0₅8345"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").Value = "'65.77"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "26.508.74"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D19").Value = "'4.738"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").Value = "'194.29"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("D21").Value = "'10.29"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").Value = "'6.353"
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'142.55"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").Value = "'0.1290"
$ws.Range("E25").Value = "  +6.37%  "
$ws.Range("D26").Value = "'7.394"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'16.28"
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("D29").Value = "'0.06334"
$ws.Range("E29").Value = "  +6.38%  "
$ws.Range("D30").Value = "'1.274"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'3.636"
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("D32").Value = "'3.457"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").Value = "'1.681"
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("D34").Value = "'1.011"
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("D35").Value = "'0.6206"
$ws.Range("E35").Value = "  +8.73%  "
$ws.Range("D36").Value = "'2.422"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").Value = "'2.789"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "'6.164"
$ws.Range("E38").Value = "  +7.51%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").Value = "1.089.01"
$ws.Range("E40").Value = "  +5.33%  "
$ws.Range("D41").Value = "'0.8652"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").Value = "'1.0000"
$ws.Range("D43").Value = "'100.43"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "1.822.27"
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("D45").Value = "'57.47"
$ws.Range("E45").Value = "  +3.01%  "
$ws.Range("D46").Value = "'8.228"
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").Value = "'0.05210"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D50").Value = "'1.486"
$ws.Range("E50").Value = "  +7.38%  "
$ws.Range("D51").Value = "'6.044"
$ws.Range("E51").Value = "  +2.35%  "

foreach ($addr in @("D5","D6","D8","D11","D14","D16","D19","D20","D21","D22","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D41","D42","D43","D45","D46","D47","D49","D50","D51")) {
    $ws.Range($addr).Style = "Normal"
}
